# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.894.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -4.00%  '
$ws.Range("E2").Style = "Normal"

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.630.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -6.58%  '
$ws.Range("E3").Style = "Normal"

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9972'
$ws.Range("D4").Style = "Normal"

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.03%  '
$ws.Range("E5").Style = "Normal"

# Row 6: USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E6").Style = "Normal"

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4725'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -6.32%  '
$ws.Range("E7").Style = "Normal"

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2543'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -7.20%  '
$ws.Range("E8").Style = "Normal"

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06088'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.59%  '
$ws.Range("E9").Style = "Normal"

# Row 10: TRON
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06964'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("E10").Style = "Normal"

# Row 11: WrappedEther
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.634.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -6.39%  '
$ws.Range("E11").Style = "Normal"

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.69'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("E12").Style = "Normal"

# Row 13: Polygon
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6117'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.52%  '
$ws.Range("E13").Style = "Normal"

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.333'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.62%  '
$ws.Range("E14").Style = "Normal"

# Row 15: Litecoin
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '72.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.75%  '
$ws.Range("E15").Style = "Normal"

# Row 17: BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9978'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E17").Style = "Normal"

# Row 18: WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.897.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.09%  '
$ws.Range("E18").Style = "Normal"

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006543'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.07%  '
$ws.Range("E19").Style = "Normal"

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -6.65%  '
$ws.Range("E20").Style = "Normal"

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.843.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.55%  '
$ws.Range("E21").Style = "Normal"

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.322'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("E22").Style = "Normal"

# Row 23: Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.526'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.90%  '
$ws.Range("E23").Style = "Normal"

# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.228'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("E24").Style = "Normal"

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E25").Style = "Normal"

# Row 26: EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.60%  '
$ws.Range("E26").Style = "Normal"

# Row 27: Toncoin
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.368'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -8.86%  '
$ws.Range("E27").Style = "Normal"

# Row 28: BitcoinCash
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '102.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("E28").Style = "Normal"

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.631'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.27%  '
$ws.Range("E29").Style = "Normal"

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.749'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("E30").Style = "Normal"

# Row 31: Stellar
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.32%  '
$ws.Range("E31").Style = "Normal"

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.526'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("E32").Style = "Normal"

# Row 33: Frax
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9992'
$ws.Range("D33").Style = "Normal"

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04267'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -8.73%  '
$ws.Range("E34").Style = "Normal"

# Row 35: HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.597'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("E35").Style = "Normal"

# Row 36: ARBITRUM
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9169'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -7.72%  '
$ws.Range("E36").Style = "Normal"

# Row 37: ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5758'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -6.96%  '
$ws.Range("E37").Style = "Normal"

# Row 38: MXToken
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -8.08%  '
$ws.Range("E38").Style = "Normal"

# Row 39: VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01536'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.68%  '
$ws.Range("E39").Style = "Normal"

# Row 40: PaxDollar
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E40").Style = "Normal"

# Row 41: TrustWalletToken
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.94%  '
$ws.Range("E41").Style = "Normal"

# Row 42: Quant
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Quant'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.09'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.92%  '
$ws.Range("E42").Style = "Normal"

# Row 43: RenderToken
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.779'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.84%  '
$ws.Range("E43").Style = "Normal"

# Row 44: TheSandbox
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3682'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.16%  '
$ws.Range("E44").Style = "Normal"

# Row 45: FraxShare
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.695'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.26%  '
$ws.Range("E45").Style = "Normal"

# Row 46: Algorand
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1088'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.92%  '
$ws.Range("E46").Style = "Normal"

# Row 47: Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05196'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("E47").Style = "Normal"

# Row 48: Aptos
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Aptos'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.032'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("E48").Style = "Normal"

# Row 49: Elrond
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Elrond'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("E49").Style = "Normal"

# Row 50: TrueUSD
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9997'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("E50").Style = "Normal"

# Row 51: USDD
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'USDD'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.23%  '
$ws.Range("E51").Style = "Normal"

